# added 4wk low sales check
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---

# MyForecast: 4-week low sales check zeroed out row 8 (W16)
$wsForecast.Range("D8").Value = 0

# Inventory Coverage (H) recalculated
$wsForecast.Range("H2").Value = 97.5
$wsForecast.Range("H3").Value = 128.67
$wsForecast.Range("H4").Value = 127.67
$wsForecast.Range("H5").Value = 126.67
$wsForecast.Range("H6").Value = 125.67
$wsForecast.Range("H7").Value = 124.67
$wsForecast.Range("H8").Value = 123.67
$wsForecast.Range("H9").Value = 184
$wsForecast.Range("H10").Value = 183
$wsForecast.Range("H11").Value = 182
$wsForecast.Range("H12").Value = 181
$wsForecast.Range("H13").Value = 180
$wsForecast.Range("H14").Value = 179
$wsForecast.Range("H15").Value = 356
$wsForecast.Range("H16").Value = 355
$wsForecast.Range("H17").Value = 354

# Seasonality Index (L) recalculated
$wsForecast.Range("L2").Value = 0.89
$wsForecast.Range("L3").Value = 0.87
$wsForecast.Range("L4").Value = 0.87
$wsForecast.Range("L5").Value = 1.18
$wsForecast.Range("L6").Value = 1.07
$wsForecast.Range("L7").Value = 0.89
$wsForecast.Range("L8").Value = 1
$wsForecast.Range("L9").Value = 0.97
$wsForecast.Range("L10").Value = 0.97
$wsForecast.Range("L11").Value = 1.16
$wsForecast.Range("L12").Value = 0.86
$wsForecast.Range("L13").Value = 0.97
$wsForecast.Range("L14").Value = 0.86
$wsForecast.Range("L15").Value = 0.92
$wsForecast.Range("L16").Value = 1.17
$wsForecast.Range("L17").Value = 0.84

# --- Summary sheet ---
$wsSummary.Range("B9").Value = "7"
$wsSummary.Range("B10").Value = "5"
$wsSummary.Range("B14").Value = "0"
